$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row of raw/clean SSA data for 2020-09-02
$row = 95
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2020-09-02"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 610957
$ws.Cells.Item($row, 3).Value = 683438
$ws.Cells.Item($row, 4).Value = 81775
$ws.Cells.Item($row, 5).Value = 65816
$ws.Cells.Item($row, 6).Value = 25.32
